$wb = $excel.ActiveWorkbook

# ALC row 104
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(104, 8).Value = 735.25
$ws.Cells.Item(104, 9).Value = 325
$ws.Cells.Item(104, 10).Value = 1145.5
$ws.Cells.Item(104, 11).Value = 975
$ws.Cells.Item(104, 12).Value = 3436.5
$ws.Cells.Item(104, 13).Value = 772
$ws.Cells.Item(104, 14).Value = -6930.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 2224.8262
$ws.Cells.Item(138, 10).Value = 2854.56
$ws.Cells.Item(138, 12).Value = 8563.68
$ws.Cells.Item(138, 14).Value = -18843.68

# ALC row 139
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(139, 8).Value = 97702.60000000001
$ws.Cells.Item(139, 10).Value = 97702.60000000001
$ws.Cells.Item(139, 12).Value = 97702.60000000001
$ws.Cells.Item(139, 14).Value = -107982.6

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(141, 8).Value = 32349.768
$ws.Cells.Item(141, 9).Value = 32349.768
$ws.Cells.Item(141, 11).Value = 97049.304
$ws.Cells.Item(141, 13).Value = -91869.304

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3699.75
$ws.Cells.Item(63, 9).Value = 2933.3333
$ws.Cells.Item(63, 10).Value = 5999
$ws.Cells.Item(63, 11).Value = 2933.3333
$ws.Cells.Item(63, 12).Value = 5999
$ws.Cells.Item(63, 13).Value = -2247.3333
$ws.Cells.Item(63, 14).Value = -7371

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 3699.75
$ws.Cells.Item(66, 9).Value = 2933.3333
$ws.Cells.Item(66, 10).Value = 5999
$ws.Cells.Item(66, 11).Value = 14666.6665
$ws.Cells.Item(66, 12).Value = 29995
$ws.Cells.Item(66, 13).Value = -11234.6665
$ws.Cells.Item(66, 14).Value = -36859

# ARM row 128
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(128, 8).Value = 77025.336
$ws.Cells.Item(128, 10).Value = 77025.336
$ws.Cells.Item(128, 12).Value = 77025.336
$ws.Cells.Item(128, 14).Value = -86985.336

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2707.9678
$ws.Cells.Item(86, 9).Value = 1786.75
$ws.Cells.Item(86, 11).Value = 1786.75
$ws.Cells.Item(86, 13).Value = -663.75

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(89, 8).Value = 2707.9678
$ws.Cells.Item(89, 9).Value = 1786.75
$ws.Cells.Item(89, 11).Value = 8933.75
$ws.Cells.Item(89, 13).Value = -3317.75

# BSM row 92
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1595.7368
$ws.Cells.Item(107, 9).Value = 1707.4375
$ws.Cells.Item(107, 11).Value = 1707.4375
$ws.Cells.Item(107, 13).Value = 212.5625

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2532.1052
$ws.Cells.Item(31, 9).Value = 1506.1875
$ws.Cells.Item(31, 10).Value = 8003.6665
$ws.Cells.Item(31, 11).Value = 1506.1875
$ws.Cells.Item(31, 12).Value = 8003.6665
$ws.Cells.Item(31, 13).Value = -1211.1875
$ws.Cells.Item(31, 14).Value = -8593.666499999999

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2532.1052
$ws.Cells.Item(34, 9).Value = 1506.1875
$ws.Cells.Item(34, 10).Value = 8003.6665
$ws.Cells.Item(34, 11).Value = 1506.1875
$ws.Cells.Item(34, 12).Value = 8003.6665
$ws.Cells.Item(34, 13).Value = -1304.1875
$ws.Cells.Item(34, 14).Value = -8407.666499999999

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(74, 8).Value = 51662.668
$ws.Cells.Item(74, 10).Value = 67494.5
$ws.Cells.Item(74, 12).Value = 67494.5
$ws.Cells.Item(74, 14).Value = -69242.5

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(77, 8).Value = 51662.668
$ws.Cells.Item(77, 10).Value = 67494.5
$ws.Cells.Item(77, 12).Value = 202483.5
$ws.Cells.Item(77, 14).Value = -211219.5

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 750
$ws.Cells.Item(22, 9).Value = 500
$ws.Cells.Item(22, 10).Value = 1000
$ws.Cells.Item(22, 11).Value = 1500
$ws.Cells.Item(22, 12).Value = 3000
$ws.Cells.Item(22, 13).Value = -1331
$ws.Cells.Item(22, 14).Value = -3338

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(27, 8).Value = 750
$ws.Cells.Item(27, 9).Value = 500
$ws.Cells.Item(27, 10).Value = 1000
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 12).Value = 3000
$ws.Cells.Item(27, 13).Value = -1398
$ws.Cells.Item(27, 14).Value = -3204

# CUL row 52
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 974.6667
$ws.Cells.Item(52, 10).Value = 974.6667
$ws.Cells.Item(52, 12).Value = 2924.0001
$ws.Cells.Item(52, 14).Value = -3456.0001

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 2664.6667
$ws.Cells.Item(122, 9).Value = 2997
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 26973
$ws.Cells.Item(122, 12).Value = 18000
$ws.Cells.Item(122, 13).Value = -24523
$ws.Cells.Item(122, 14).Value = -22900

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4930.6
$ws.Cells.Item(70, 9).Value = 4254
$ws.Cells.Item(70, 10).Value = 5099.75
$ws.Cells.Item(70, 11).Value = 4254
$ws.Cells.Item(70, 12).Value = 5099.75
$ws.Cells.Item(70, 13).Value = -3984
$ws.Cells.Item(70, 14).Value = -5639.75

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 4930.6
$ws.Cells.Item(73, 9).Value = 4254
$ws.Cells.Item(73, 10).Value = 5099.75
$ws.Cells.Item(73, 11).Value = 4254
$ws.Cells.Item(73, 12).Value = 5099.75
$ws.Cells.Item(73, 13).Value = -3318
$ws.Cells.Item(73, 14).Value = -6971.75

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2401.8096
$ws.Cells.Item(132, 9).Value = 2233.4119
$ws.Cells.Item(132, 10).Value = 3117.5
$ws.Cells.Item(132, 11).Value = 6700.2357
$ws.Cells.Item(132, 12).Value = 9352.5
$ws.Cells.Item(132, 13).Value = -4170.2357
$ws.Cells.Item(132, 14).Value = -14412.5

# LTW row 12
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12, 8).Value = 5000
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 5000
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 5000
$ws.Cells.Item(12, 14).Value = -5340
$ws.Cells.Item(12, 13).ClearContents()

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 948.2273
$ws.Cells.Item(82, 9).Value = 956.82355
$ws.Cells.Item(82, 11).Value = 956.82355
$ws.Cells.Item(82, 13).Value = -595.82355

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(85, 8).Value = 948.2273
$ws.Cells.Item(85, 9).Value = 956.82355
$ws.Cells.Item(85, 11).Value = 956.82355
$ws.Cells.Item(85, 13).Value = 291.17645

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 25166.928
$ws.Cells.Item(93, 9).Value = 942.6
$ws.Cells.Item(93, 11).Value = 942.6
$ws.Cells.Item(93, 13).Value = 305.4

# LTW row 107
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(107, 8).Value = 4999
$ws.Cells.Item(107, 9).Value = 4999
$ws.Cells.Item(107, 11).Value = 4999
$ws.Cells.Item(107, 13).Value = -3079

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 3806
$ws.Cells.Item(122, 9).Value = 3590.8333
$ws.Cells.Item(122, 11).Value = 10772.4999
$ws.Cells.Item(122, 13).Value = -8322.499899999999

# WVR row 23
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 2999.6667
$ws.Cells.Item(23, 9).Value = 2999.5
$ws.Cells.Item(23, 10).Value = 3000
$ws.Cells.Item(23, 11).Value = 2999.5
$ws.Cells.Item(23, 12).Value = 3000
$ws.Cells.Item(23, 13).Value = -2770.5
$ws.Cells.Item(23, 14).Value = -3458

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3671.8845
$ws.Cells.Item(81, 9).Value = 3698.76
$ws.Cells.Item(81, 11).Value = 7397.52
$ws.Cells.Item(81, 13).Value = -6336.52

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 3671.8845
$ws.Cells.Item(84, 9).Value = 3698.76
$ws.Cells.Item(84, 11).Value = 36987.60000000001
$ws.Cells.Item(84, 13).Value = -31683.60000000001

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 3128
$ws.Cells.Item(100, 9).Value = 3339.6365
$ws.Cells.Item(100, 11).Value = 6679.273
$ws.Cells.Item(100, 13).Value = -6138.273

# WVR row 124
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(124, 8).Value = 124460.5
$ws.Cells.Item(124, 10).Value = 124460.5
$ws.Cells.Item(124, 12).Value = 124460.5
$ws.Cells.Item(124, 14).Value = -134280.5

# WVR row 125
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(125, 8).Value = 76666.664
$ws.Cells.Item(125, 10).Value = 76666.664
$ws.Cells.Item(125, 12).Value = 76666.664
$ws.Cells.Item(125, 14).Value = -86506.664

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3500.476
$ws.Cells.Item(126, 9).Value = 2861.9443
$ws.Cells.Item(126, 10).Value = 7331.6665
$ws.Cells.Item(126, 11).Value = 8585.832900000001
$ws.Cells.Item(126, 12).Value = 21994.9995
$ws.Cells.Item(126, 13).Value = -6115.832900000001
$ws.Cells.Item(126, 14).Value = -26934.9995

# WVR row 129
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(129, 8).Value = 98136.664
$ws.Cells.Item(129, 10).Value = 98705
$ws.Cells.Item(129, 12).Value = 98705
$ws.Cells.Item(129, 14).Value = -108705

# WVR row 130
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(130, 8).Value = 74849
$ws.Cells.Item(130, 10).Value = 74849
$ws.Cells.Item(130, 12).Value = 74849
$ws.Cells.Item(130, 14).Value = -84889

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 7618.9165
$ws.Cells.Item(132, 9).Value = 4553.9697
$ws.Cells.Item(132, 11).Value = 13661.9091
$ws.Cells.Item(132, 13).Value = -11131.9091

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(135, 8).Value = 97857.5
$ws.Cells.Item(135, 10).Value = 97857.5
$ws.Cells.Item(135, 12).Value = 97857.5
$ws.Cells.Item(135, 14).Value = -107997.5

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 328.9
$ws.Cells.Item(136, 9).Value = 328.9
$ws.Cells.Item(136, 11).Value = 986.6999999999999
$ws.Cells.Item(136, 13).Value = 1563.3
